# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 44 and 45) at the top of the data
# block for "Frutilla" @ Vega Modelo de Temuco, pushing the existing
# rows (old 44-128) down by two (new 46-130).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 44.
$ws.Rows("44:45").Insert()

# --- New row 44 ---
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = 44469
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100101
$ws.Range("H44").Value = "Berries"
$ws.Range("I44").Value = 100112025
$ws.Range("J44").Value = "Frutilla"
$ws.Range("K44").Value = "Sin especificar"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 3000
$ws.Range("N44").Value = 14000
$ws.Range("O44").Value = 16000
$ws.Range("P44").Value = 15333
$ws.Range("Q44").Value = "$/bandeja 7 kilos"
$ws.Range("R44").Value = "Provincia de Melipilla"
$ws.Range("S44").Value = 2190
$ws.Range("T44").Value = 7

# --- New row 45 ---
$ws.Range("A45").Value = 10
$ws.Range("B45").Value = "Vega Modelo de Temuco"
$ws.Range("C45").Value = "La Araucanía"
$ws.Range("D45").Value = 44469
$ws.Range("E45").Value = 9
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100101
$ws.Range("H45").Value = "Berries"
$ws.Range("I45").Value = 100112025
$ws.Range("J45").Value = "Frutilla"
$ws.Range("K45").Value = "Sin especificar"
$ws.Range("L45").Value = "Tercera"
$ws.Range("M45").Value = 100
$ws.Range("N45").Value = 7000
$ws.Range("O45").Value = 7000
$ws.Range("P45").Value = 7000
$ws.Range("Q45").Value = "$/bandeja 7 kilos"
$ws.Range("R45").Value = "Provincia de Melipilla"
$ws.Range("S45").Value = 1000
$ws.Range("T45").Value = 7
